$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 246, shifting existing rows 246:260 down to 247:261
$ws.Rows.Item(246).Insert()

# Populate the newly inserted row 246 with the new record
$ws.Cells.Item(246, 1).Value = 10
$ws.Cells.Item(246, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(246, 3).Value = "La Araucanía"
$ws.Cells.Item(246, 4).Value = 44585
$ws.Cells.Item(246, 5).Value = 9
$ws.Cells.Item(246, 6).Value = 100112044
$ws.Cells.Item(246, 7).Value = "Perejil"
$ws.Cells.Item(246, 8).Value = "Sin especificar"
$ws.Cells.Item(246, 9).Value = "Primera"
$ws.Cells.Item(246, 10).Value = 65
$ws.Cells.Item(246, 11).Value = 5000
$ws.Cells.Item(246, 12).Value = 5000
$ws.Cells.Item(246, 13).Value = 5000
$ws.Cells.Item(246, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(246, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(246, 16).Value = 1667
$ws.Cells.Item(246, 17).Value = 3
$ws.Cells.Item(246, 18).Value = "Hortaliza"
